# Generate Report for Handback
# Updates the handback-status workbook with a new pair of handed-back files
# (new GUIDs / new content hash / new handoff+handback timestamps) in place
# of the previous pair, across the Overview, zh-cn and de-de sheets -
# including both the cell text and the matching hyperlink display text.

$wb = $excel.ActiveWorkbook

# Map of old literal cell/hyperlink text -> new literal text.
$replacements = @{
    "8c50de45-616d-4b0f-9a5b-ab47a1647522.md" = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.md";
    "f28a17c2-0989-40c6-852a-541543414afe.md" = "ffffc0ee9c18-8472-429c-b9df-d7d428c056a5.md";

    "8c50de45-616d-4b0f-9a5b-ab47a1647522.61a71177ba8f505012c78f20b07b71b5dd019a8e.zh-cn.xlf" = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.zh-cn.xlf";
    "f28a17c2-0989-40c6-852a-541543414afe.99985fb90f70a4870b2223f4a87f5a18236a1af3.zh-cn.xlf" = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.zh-cn.xlf";

    "8c50de45-616d-4b0f-9a5b-ab47a1647522.61a71177ba8f505012c78f20b07b71b5dd019a8e.de-de.xlf" = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.de-de.xlf";
    "f28a17c2-0989-40c6-852a-541543414afe.99985fb90f70a4870b2223f4a87f5a18236a1af3.de-de.xlf" = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.de-de.xlf";

    "2016-03-12 14:47:47" = "2016-03-12 14:48:51";
    "2016-03-12 14:48:05" = "2016-03-12 14:49:08";
    "2016-03-12 14:47:50" = "2016-03-12 14:48:56";
    "2016-03-12 14:48:11" = "2016-03-12 14:49:15";
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            $val = $cell.Value2
            if ($null -ne $val -and $replacements.ContainsKey($val)) {
                $cell.Value2 = $replacements[$val]
            }
        }
    }

    foreach ($hl in $ws.Hyperlinks) {
        $disp = $hl.TextToDisplay
        if ($null -ne $disp -and $replacements.ContainsKey($disp)) {
            $hl.TextToDisplay = $replacements[$disp]
        }
    }
}
